$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("C18").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("D18").Value = "[-, 'ELM-2NA-Instalções Elétricas']"
$ws.Range("E18").Value = "[-, 'ELM-2NA-Instalções Elétricas']"
$ws.Range("F18").Value = "['ELM-2NA-Lab. De Máquinas elétricas', 'ELM-2NA-Lab. De Máquinas elétricas']"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("C19").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("E19").Value = "-"

# Row 20
$ws.Range("E20").Value = "['ELM-2NA-Instalções Elétricas', -]"
$ws.Range("F20").Value = "ELM-1NA-Circuitos Elétricos 1"

# Row 21
$ws.Range("B21").Value = "['ELM-2NA-Instalções Elétricas', -]"
$ws.Range("F21").Value = "ELM-1NA-Circuitos Elétricos 1"
